$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the missing date + hours entry for row 7 (continues the Date/Hours
# series started in rows 5-6), which feeds the D5 SUM and F5 totals.
$ws.Range("A7").Value = 41206
$ws.Range("B7").Value = 6
